# Rename landings_observed/landings_predicted to the *_weight variants, and
# add the new catch_weight / catch_numbers / discards_weight labels that go
# alongside the existing CATCH rows (ret_bio / ret_num / discard).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CATCH / obs -> landings_observed_weight (was landings_observed)
$ws.Range("C59").Value = "landings_observed_weight"

# CATCH / exp -> landings_predicted_weight (was landings_predicted)
$ws.Range("C61").Value = "landings_predicted_weight"

# CATCH / ret_bio gains an alt_label: catch_weight
$ws.Range("C70").Value = "catch_weight"

# CATCH / ret_num gains an alt_label: catch_numbers
$ws.Range("C73").Value = "catch_numbers"

# discard gains an alt_label: discards_weight
$ws.Range("C74").Value = "discards_weight"

# Update the view's active cell / selection to match the author's edit.
$ws.Range("E61").Select()
